$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold numeric-looking values stored as text (shared
# strings), e.g. "217.5". Assigning a plain numeric-looking string would
# make Excel auto-convert the cell to a real number, so we temporarily
# force a text number format, assign the new text values, then clear the
# formatting again so the cells end up with no explicit style applied -
# matching the original workbook where these cells carry no style index.

$ws.Range("B2:B9").NumberFormat = "@"

$ws.Range("B2").Value = "205.9"
$ws.Range("B3").Value = "200.7"
$ws.Range("B4").Value = "211.1"
$ws.Range("B5").Value = "207.9"
$ws.Range("B6").Value = "186.7"
$ws.Range("B7").Value = "182.7"
$ws.Range("B8").Value = "190.4"
$ws.Range("B9").Value = "188.3"

$ws.Range("B2:B9").ClearFormats()
